# Fix agent-column documentation typos / inconsistencies in the "Data" sheet
# of the ArchivesSpace import Excel template.
#
# Summary of the edit:
#  - The 3rd "Person Agent" block (columns AS/AT) had "Role" and
#    "header string" swapped relative to the 1st/2nd Person Agent blocks.
#    Swap the contents (and column widths) of columns AS and AT on rows 3-5
#    so Agent(3) matches the order used by Agent(1)/Agent(2).
#  - The "Corporate Agent" (1) block (columns AZ:BC) is relabeled with an
#    explicit "(1)" suffix to match the "(2)" suffix already used by the
#    following Corporate Agent (2) block (columns BD:BG).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Swap columns AS (45) and AT (46) for rows 3, 4 and 5 -------------------
foreach ($row in 3..5) {
    $asCell = $ws.Cells.Item($row, 45)
    $atCell = $ws.Cells.Item($row, 46)
    $tmp = $asCell.Value
    $asCell.Value = $atCell.Value
    $atCell.Value = $tmp
}

# Swap the column widths too, so AS/AT keep matching their (swapped) content
$asWidth = $ws.Columns.Item(45).ColumnWidth
$atWidth = $ws.Columns.Item(46).ColumnWidth
$ws.Columns.Item(45).ColumnWidth = $atWidth
$ws.Columns.Item(46).ColumnWidth = $asWidth

# --- Relabel the Corporate Agent (1) block (AZ5:BB5) with "(1)" suffixes ----
$ws.Cells.Item(5, 52).Value = "Corporate Agent Record (1) ID"
$ws.Cells.Item(5, 53).Value = "Corporate Agent(1) header string"
$ws.Cells.Item(5, 54).Value = "Corporate Agent(1) Role"

# --- Update the "Role" data validation list to follow column AT now --------
$roleValidationRange = $ws.Range("AL6:AL1048576,AP6:AP1048576,AX6:AX1048576,BB6:BB1048576,BF6:BF1048576,AT6:AT1048576")
$roleValidationRange.Validation.Delete()
$roleValidationRange.Validation.Add(3, 1, 1, "=Sheet2!`$A`$2:`$A`$4")

# --- Reflect the selection/view state after performing the column swap -----
$ws.Range("AS1:AS1048576").Select()
$ws.Application.ActiveWindow.ScrollColumn = 36
